$d = $word.ActiveDocument

# Helper: given a Range positioned at an (empty) insertion point, insert a
# sequence of text fragments as separate runs (one w:r per fragment),
# advancing the range after each insert so the next fragment lands after
# the previous one.
function Insert-Runs($rng, $fragments) {
    foreach ($frag in $fragments) {
        $rng.InsertAfter($frag)
        $rng.Collapse(0)
    }
}

# --- Paragraph 1 (GOAL -> NOTE rule) ---------------------------------
# " create any number of GOALs. Each GOAL may only be created by one and
#   only one USER." becomes a run-split rewrite about NOTEs.
$rng = $d.Content
$rng.Find.Execute("create any number of GOALs. Each GOAL may only be created by one and only one USER.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = ""
$rng.Collapse(0)
Insert-Runs $rng @(
    " ",
    "write",
    " any number of ",
    "NOTEs",
    ". Each ",
    "NOTE",
    " may only be created ",
    "and edited ",
    "by one and only one USER."
)

# --- Paragraph 2 (NOTE -> EXCERCISE rule) ----------------------------
$rng2 = $d.Content
$rng2.Find.Execute("A USER may write any number of NOTEs. Each NOTE may be created by only one USER.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Text = ""
$rng2.Collapse(0)
Insert-Runs $rng2 @(
    "A USER may ",
    "create",
    " ",
    "up to 5",
    " ",
    "EXCERCISEs",
    ". ",
    "Each ",
    "EXCERCISE",
    " may ",
    "be viewed ",
    "and created ",
    "by one and only one USER."
)

# --- Paragraph 3 (BUDGET rule) ---------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("A USER may or may not submit one BUDGET. A BUDGET may be submitted by only one USER>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3.Text = ""
$rng3.Collapse(0)
Insert-Runs $rng3 @(
    "A USER may or may not submit one BUDGET. A BUDGET may be ",
    "submitted",
    " by ",
    "one and ",
    "only one USER",
    "."
)
